$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

$ws.Range("E5").Value = 169
$ws.Range("F5").Value = 115
$ws.Range("H5").Value = 126

$ws.Range("E9").Value = 14

$ws.Range("E10").Value = 712
$ws.Range("F10").Value = 396
$ws.Range("H10").Value = 491

$ws.Range("E11").Value = 472
$ws.Range("F11").Value = 265
$ws.Range("H11").Value = 330

$ws.Range("E12").Value = 711
$ws.Range("F12").Value = 425
$ws.Range("H12").Value = 511

$ws.Range("E13").Value = 169

$ws.Range("E14").Value = 144

$ws.Range("F15").Value = 97
$ws.Range("H15").Value = 148

$ws.Range("E16").Value = 234
$ws.Range("F16").Value = 134
$ws.Range("H16").Value = 182

$ws.Range("E17").Value = 126
$ws.Range("F17").Value = 69
$ws.Range("H17").Value = 93

$ws.Range("E22").Value = 197

$ws.Range("E23").Value = 229
$ws.Range("F23").Value = 117
$ws.Range("H23").Value = 169

$ws.Range("E24").Value = 272
$ws.Range("F24").Value = 159
$ws.Range("H24").Value = 189

$ws.Range("E25").Value = 336
$ws.Range("F25").Value = 185
$ws.Range("H25").Value = 245

$ws.Range("E26").Value = 203

$ws.Range("E27").Value = 386

$ws.Range("E28").Value = 230

$ws.Range("E30").Value = 256

$ws.Range("E32").Value = 217
$ws.Range("F32").Value = 137
$ws.Range("H32").Value = 175

$ws.Range("E33").Value = 335
$ws.Range("G33").Value = 91
$ws.Range("H33").Value = 271

$ws.Range("E34").Value = 256
$ws.Range("F34").Value = 180
$ws.Range("H34").Value = 218

$ws.Range("E35").Value = 187

$ws.Range("G37").Value = 36
$ws.Range("H37").Value = 147

$ws.Range("F38").Value = 65
$ws.Range("H38").Value = 82

$ws.Range("E39").Value = 202
$ws.Range("F39").Value = 105
$ws.Range("H39").Value = 156

$ws.Range("E40").Value = 311
$ws.Range("F40").Value = 159
$ws.Range("H40").Value = 239

$ws.Range("F41").Value = 224
$ws.Range("H41").Value = 316

$ws.Range("E42").Value = 468
$ws.Range("F42").Value = 265
$ws.Range("H42").Value = 326

$ws.Range("E43").Value = 145
$ws.Range("F43").Value = 83
$ws.Range("H43").Value = 110

$ws.Range("E44").Value = 379
$ws.Range("F44").Value = 196
$ws.Range("H44").Value = 264

$ws.Range("E45").Value = 184

$ws.Range("E46").Value = 395
$ws.Range("F46").Value = 231
$ws.Range("H46").Value = 295

$ws.Range("E48").Value = 276
$ws.Range("F48").Value = 134
$ws.Range("H48").Value = 178

$ws.Range("E49").Value = 342
